$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the email address string in A2 with a plain numeric value
$ws.Range("A2").Value = 12345.258

# Move/leave the active selection on A2 (was A9)
$ws.Range("A2").Select()
